$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.670167666666667
$ws.Range("H2").Value = 20.010503
$ws.Range("I2").Value = 0.0423069620011633
$ws.Range("J2").Value = 0.0423069620011633
$ws.Range("M2").Value = 0.11549
$ws.Range("N2").Value = 0.34647
$ws.Range("O2").Value = 0.01449407350231777
$ws.Range("P2").Value = 0.01449407350231777
$ws.Range("Q2").Value = 0.7703376638233334
$ws.Range("R2").Value = 6.93303897441
$ws.Range("S2").Value = 0.0006132002169046255
$ws.Range("T2").Value = 0.0006132002169046256
$ws.Range("G3").Value = 6.670167666666667
$ws.Range("H3").Value = 20.010503
$ws.Range("I3").Value = 0.0423069620011633
$ws.Range("J3").Value = 0.0423069620011633
$ws.Range("O3").Value = 0.9654500393716549
$ws.Range("P3").Value = 0.965450039371655
$ws.Range("Q3").Value = 51.3121813373429
$ws.Range("R3").Value = 461.809632036086
$ws.Range("S3").Value = 0.04084525812971821
$ws.Range("T3").Value = 0.04084525812971821
$ws.Range("G4").Value = 6.670167666666667
$ws.Range("H4").Value = 20.010503
$ws.Range("I4").Value = 0.0423069620011633
$ws.Range("J4").Value = 0.0423069620011633
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.159807
$ws.Range("N4").Value = 0.479421
$ws.Range("O4").Value = 0.02005588712602732
$ws.Range("P4").Value = 0.02005588712602732
$ws.Range("Q4").Value = 1.065939484307
$ws.Range("R4").Value = 9.593455358763
$ws.Range("S4").Value = 0.000848503654540458
$ws.Range("T4").Value = 0.000848503654540458
$ws.Range("I5").Value = 0.9513278459982415
$ws.Range("J5").Value = 0.9513278459982416
$ws.Range("M5").Value = 0.11549
$ws.Range("N5").Value = 0.34647
$ws.Range("O5").Value = 0.01449407350231777
$ws.Range("P5").Value = 0.01449407350231777
$ws.Range("Q5").Value = 17.32205849231667
$ws.Range("R5").Value = 155.89852643085
$ws.Range("S5").Value = 0.01378861572470015
$ws.Range("T5").Value = 0.01378861572470015
$ws.Range("I6").Value = 0.9513278459982415
$ws.Range("J6").Value = 0.9513278459982416
$ws.Range("O6").Value = 0.9654500393716549
$ws.Range("P6").Value = 0.965450039371655
$ws.Range("S6").Value = 0.9184595063743538
$ws.Range("T6").Value = 0.918459506374354
$ws.Range("I7").Value = 0.9513278459982415
$ws.Range("J7").Value = 0.9513278459982416
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.159807
$ws.Range("N7").Value = 0.479421
$ws.Range("O7").Value = 0.02005588712602732
$ws.Range("P7").Value = 0.02005588712602732
$ws.Range("Q7").Value = 23.969055342295
$ws.Range("R7").Value = 215.721498080655
$ws.Range("S7").Value = 0.01907972389918743
$ws.Range("T7").Value = 0.01907972389918744
$ws.Range("G8").Value = 1.003544
$ws.Range("H8").Value = 3.010632
$ws.Range("I8").Value = 0.0063651920005952
$ws.Range("J8").Value = 0.0063651920005952
$ws.Range("M8").Value = 0.11549
$ws.Range("N8").Value = 0.34647
$ws.Range("O8").Value = 0.01449407350231777
$ws.Range("P8").Value = 0.01449407350231777
$ws.Range("Q8").Value = 0.11589929656
$ws.Range("R8").Value = 1.04309366904
$ws.Range("S8").Value = 0.00009225756071299189
$ws.Range("T8").Value = 0.00009225756071299191
$ws.Range("G9").Value = 1.003544
$ws.Range("H9").Value = 3.010632
$ws.Range("I9").Value = 0.0063651920005952
$ws.Range("J9").Value = 0.0063651920005952
$ws.Range("O9").Value = 0.9654500393716549
$ws.Range("P9").Value = 0.965450039371655
$ws.Range("Q9").Value = 7.720050571642668
$ws.Range("R9").Value = 69.48045514478402
$ws.Range("S9").Value = 0.006145274867582778
$ws.Range("T9").Value = 0.00614527486758278
$ws.Range("G10").Value = 1.003544
$ws.Range("H10").Value = 3.010632
$ws.Range("I10").Value = 0.0063651920005952
$ws.Range("J10").Value = 0.0063651920005952
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.159807
$ws.Range("N10").Value = 0.479421
$ws.Range("O10").Value = 0.02005588712602732
$ws.Range("P10").Value = 0.02005588712602732
$ws.Range("Q10").Value = 0.160373356008
$ws.Range("R10").Value = 1.443360204072
$ws.Range("S10").Value = 0.0001276595722994294
$ws.Range("T10").Value = 0.0001276595722994294
